$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the backlog items currently in rows 5-9: the "Interaktiv sidlayout /
# Produktsida / Koplayout" (sprint S0) trio moves above the "Skapa Produktstruct /
# Koppla Databas till Sida" pair (sprint S1), which now follow them.
$ws.Range("B5").Value = "Interaktiv sidlayout"
$ws.Range("C5").Value = "Front end"
$ws.Range("D5").Value = "Låg"
$ws.Range("F5").Value = "S0"

$ws.Range("B6").Value = "Produktsida"
$ws.Range("C6").Value = "Front end"
$ws.Range("D6").Value = "Medel"
$ws.Range("F6").Value = "S0"

$ws.Range("B7").Value = "Köplayout"
$ws.Range("C7").Value = "Front end"
$ws.Range("D7").Value = "Låg"
$ws.Range("F7").Value = "S0"

$ws.Range("B8").Value = "Skapa Produktstruct"
$ws.Range("C8").Value = "Databas"
$ws.Range("D8").Value = "Hög"
$ws.Range("F8").Value = "S1"

$ws.Range("B9").Value = "Koppla Databas till Sida"
$ws.Range("C9").Value = "Databas"
$ws.Range("D9").Value = "Hög"
$ws.Range("F9").Value = "S1"

# Drop the stray "kolla /index.php" comment next to the PHP-script task.
$ws.Range("G10").ClearContents()

# Shoppingvagn/Recension-layout and placeholder-swap tasks are now underway
# (sprint S1) instead of merely planned.
$ws.Range("E11").Style = "Bra"
$ws.Range("F11").Value = "S1"

$ws.Range("E12").Style = "Bra"
$ws.Range("F12").Value = "S1"

# Restore the user's last selection / scroll position.
$ws.Range("G11").Select()
